$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '43.145.02'
$ws.Range("E2").Value = '  -0.62%  '
$ws.Range("D3").Value = '2.314.24'
$ws.Range("E3").Value = '  -2.46%  '
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("D5").Value = '303.64'
$ws.Range("E5").Value = '  -1.87%  '
$ws.Range("D6").Value = '99.89'
$ws.Range("E6").Value = '  -4.91%  '
$ws.Range("D7").Value = '0.509'
$ws.Range("E7").Value = '  -1.73%  '
$ws.Range("E8").Value = '  +0.07%  '
$ws.Range("D9").Value = '0.509'
$ws.Range("E9").Value = '  -1.66%  '
$ws.Range("D10").Value = '34.79'
$ws.Range("E10").Value = '  -3.88%  '
$ws.Range("D11").Value = '51.18'
$ws.Range("E11").Value = '  -4.04%  '
$ws.Range("D12").Value = '0.0795'
$ws.Range("E12").Value = '  -2.14%  '
$ws.Range("E13").Value = '  +0.41%  '
$ws.Range("D14").Value = '6.78'
$ws.Range("E14").Value = '  -2.98%  '
$ws.Range("D15").Value = '2.681.76'
$ws.Range("E15").Value = '  -2.14%  '
$ws.Range("D16").Value = '15.57'
$ws.Range("E16").Value = '  -0.19%  '
$ws.Range("D17").Value = '2.318.82'
$ws.Range("E17").Value = '  -2.26%  '
$ws.Range("D18").Value = '0.799'
$ws.Range("E18").Value = '  -1.80%  '
$ws.Range("D19").Value = '43.049.99'
$ws.Range("E19").Value = '  -0.78%  '
$ws.Range("D20").Value = '11.75'
$ws.Range("E20").Value = '  -2.38%  '
$ws.Range("D21").Value = '0.0₃0901'
$ws.Range("E21").Value = '  -1.85%  '
$ws.Range("D22").Value = '6.07'
$ws.Range("E22").Value = '  -3.68%  '
$ws.Range("D23").Value = '67.44'
$ws.Range("E23").Value = '  -1.19%  '
$ws.Range("D24").Value = '237.51'
$ws.Range("E24").Value = '  -1.67%  '
$ws.Range("D25").Value = '1.96'
$ws.Range("E25").Value = '  -4.52%  '
$ws.Range("D26").Value = '2.52'
$ws.Range("E26").Value = '  -3.85%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").Value = '24.90'
$ws.Range("E28").Value = '  -3.57%  '
$ws.Range("D29").Value = '2.17'
$ws.Range("E29").Value = '  -5.68%  '
$ws.Range("D30").Value = '34.42'
$ws.Range("E30").Value = '  -6.35%  '
$ws.Range("D31").Value = '165.19'
$ws.Range("E31").Value = '  +2.17%  '
$ws.Range("D32").Value = '9.19'
$ws.Range("E32").Value = '  -3.86%  '
$ws.Range("D33").Value = '1.00'
$ws.Range("E33").Value = '  +0.10%  '
$ws.Range("D34").Value = '5.04'
$ws.Range("E34").Value = '  -4.37%  '
$ws.Range("E35").Value = '  -4.83%  '
$ws.Range("D36").Value = '4.49'
$ws.Range("E36").Value = '  -5.80%  '
$ws.Range("B37").Value = 'Hedera'
$ws.Range("C37").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D37").Value = '0.0702'
$ws.Range("E37").Value = '  -5.49%  '
$ws.Range("B38").Value = 'Celestia'
$ws.Range("C38").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D38").Value = '16.75'
$ws.Range("E38").Value = '  -8.91%  '
$ws.Range("D39").Value = '2.88'
$ws.Range("E39").Value = '  -7.33%  '
$ws.Range("E40").Value = '  -6.83%  '
$ws.Range("E41").Value = '  -4.05%  '
$ws.Range("E42").Value = '  -3.13%  '
$ws.Range("E43").Value = '  -9.53%  '
$ws.Range("D44").Value = '1.975.23'
$ws.Range("E44").Value = '  -1.55%  '
$ws.Range("D45").Value = '0.0283'
$ws.Range("E45").Value = '  -2.45%  '
$ws.Range("D46").Value = '18.37'
$ws.Range("E46").Value = '  -7.74%  '
$ws.Range("D47").Value = '2.92'
$ws.Range("E47").Value = '  -7.10%  '
$ws.Range("D48").Value = '9.78'
$ws.Range("E48").Value = '  -7.91%  '
$ws.Range("D49").Value = '54.82'
$ws.Range("E49").Value = '  -5.71%  '
$ws.Range("D50").Value = '4.82'
$ws.Range("E50").Value = '  +2.17%  '
$ws.Range("D51").Value = '2.543.56'
$ws.Range("E51").Value = '  -0.51%  '
